# Commit: "modified born position of city"
#
# The "Scene" sheet has a header row (row 1) with columns:
#   A=ID, B=FilePath, C=MaxGroup, D=MaxGroupPlayers, E=RelivePos,
#   F=SceneName, G=SceneShowName, H=LoadingUI, I=SoundList,
#   J=CamOffestPos, K=CamOffestRot, L=Width, M=CanClone, N=ActorID
#
# Row 2 is the "villageScene" (city) scene, ID=1. Its RelivePos (born /
# respawn position) is updated from "0,0,0" to "20,0,-137".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = "20,0,-137"
